# Weekly fruit/vegetable data update: insert a new price record as row 192,
# shifting all subsequent rows down by one (old row 192 becomes 193, ...,
# old row 273 becomes 274).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 192; Excel pushes existing rows 192..273 down
# to 193..274 and copies formatting (e.g. the date format on column D) from
# the row above.
$ws.Rows("192:192").Insert()

# Populate the newly inserted row 192 with the new record.
$ws.Range("A192").Value2 = 3
$ws.Range("B192").Value2 = "Femacal de La Calera"
$ws.Range("C192").Value2 = "Coquimbo"
$ws.Range("D192").Value2 = 44466
$ws.Range("E192").Value2 = 5
$ws.Range("F192").Value2 = 100112021
$ws.Range("G192").Value2 = "Ají"
$ws.Range("H192").Value2 = "Americana (o)"
$ws.Range("I192").Value2 = "Primera"
$ws.Range("J192").Value2 = 53
$ws.Range("K192").Value2 = 43000
$ws.Range("L192").Value2 = 44000
$ws.Range("M192").Value2 = 43472
$ws.Range("N192").Value2 = "`$/caja 15 kilos"
$ws.Range("O192").Value2 = "Región de Arica y Parinacota"
$ws.Range("P192").Value2 = 2898
$ws.Range("Q192").Value2 = 15
$ws.Range("R192").Value2 = "Hortaliza"
